$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + '36.568.74'
$ws.Range("E2").Value = '  +0.20%  '
$ws.Range("D3").Value = "'" + '1.954.92'
$ws.Range("E3").Value = '  +0.74%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = "'" + '243.68'
$ws.Range("E5").Value = '  +0.27%  '
$ws.Range("D6").Value = "'" + '0.624'
$ws.Range("E6").Value = '  +2.02%  '
$ws.Range("D7").Value = "'" + '60.10'
$ws.Range("E7").Value = '  +5.27%  '
$ws.Range("E8").Value = '  -0.07%  '
$ws.Range("E9").Value = '  +4.45%  '
$ws.Range("D10").Value = "'" + '0.0790'
$ws.Range("E10").Value = '  -5.28%  '
$ws.Range("E11").Value = '  +0.50%  '
$ws.Range("D12").Value = "'" + '14.22'
$ws.Range("E12").Value = '  +6.62%  '
$ws.Range("D13").Value = "'" + '2.239.65'
$ws.Range("E13").Value = '  +0.58%  '
$ws.Range("D14").Value = "'" + '0.826'
$ws.Range("E14").Value = '  +2.31%  '
$ws.Range("D15").Value = "'" + '21.56'
$ws.Range("E15").Value = '  +1.65%  '
$ws.Range("D16").Value = "'" + '5.25'
$ws.Range("E16").Value = '  +1.92%  '
$ws.Range("D17").Value = "'" + '1.958.65'
$ws.Range("E17").Value = '  +0.87%  '
$ws.Range("D18").Value = "'" + '36.495.51'
$ws.Range("E18").Value = '  +0.20%  '
$ws.Range("D19").Value = "'" + '69.21'
$ws.Range("E19").Value = '  -0.05%  '
$ws.Range("D20").Value = "'" + '0.0₃0852'
$ws.Range("E20").Value = '  -0.74%  '
$ws.Range("D21").Value = "'" + '229.32'
$ws.Range("E21").Value = '  +0.78%  '
$ws.Range("D22").Value = "'" + '5.07'
$ws.Range("E22").Value = '  +1.88%  '
$ws.Range("E23").Value = '  +0.12%  '
$ws.Range("E24").Value = '  +0.44%  '
$ws.Range("E25").Value = '  +3.13%  '
$ws.Range("E26").Value = '  +8.57%  '
$ws.Range("D27").Value = "'" + '9.15'
$ws.Range("E27").Value = '  +0.19%  '
$ws.Range("D28").Value = "'" + '160.05'
$ws.Range("E28").Value = '  +0.32%  '
$ws.Range("D29").Value = "'" + '19.27'
$ws.Range("E29").Value = '  +0.76%  '
$ws.Range("D30").Value = "'" + '1.31'
$ws.Range("E30").Value = '  +20.93%  '
$ws.Range("E31").Value = '  +1.86%  '
$ws.Range("E32").Value = '  +3.79%  '
$ws.Range("D33").Value = "'" + '0.0611'
$ws.Range("E33").Value = '  -0.32%  '
$ws.Range("D34").Value = "'" + '4.45'
$ws.Range("E34").Value = '  +7.29%  '
$ws.Range("B35").Value = 'LidoDAOToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D35").Value = "'" + '2.27'
$ws.Range("E35").Value = '  +4.37%  '
$ws.Range("B36").Value = 'BinanceUSD'
$ws.Range("C36").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D36").Value = "'" + '1.00'
$ws.Range("E36").Value = '  -0.16%  '
$ws.Range("D37").Value = "'" + '3.43'
$ws.Range("E37").Value = '  +10.25%  '
$ws.Range("E38").Value = '  -0.66%  '
$ws.Range("D39").Value = "'" + '5.44'
$ws.Range("E39").Value = '  -11.30%  '
$ws.Range("D40").Value = "'" + '0.0966'
$ws.Range("E40").Value = '  -2.12%  '
$ws.Range("E41").Value = '  -0.29%  '
$ws.Range("E42").Value = '  +1.73%  '
$ws.Range("D43").Value = "'" + '0.0210'
$ws.Range("E43").Value = '  +0.64%  '
$ws.Range("D44").Value = "'" + '15.77'
$ws.Range("E44").Value = '  +0.45%  '
$ws.Range("D45").Value = "'" + '1.364.10'
$ws.Range("E45").Value = '  +2.08%  '
$ws.Range("D46").Value = "'" + '88.71'
$ws.Range("E46").Value = '  +3.23%  '
$ws.Range("D47").Value = "'" + '1.03'
$ws.Range("E47").Value = '  +0.00%  '
$ws.Range("D48").Value = "'" + '7.13'
$ws.Range("E48").Value = '  +0.08%  '
$ws.Range("E49").Value = '  +0.70%  '
$ws.Range("D50").Value = "'" + '45.22'
$ws.Range("E50").Value = '  +4.80%  '
$ws.Range("D51").Value = "'" + '2.135.51'
$ws.Range("E51").Value = '  +0.83%  '
